$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# The sheet currently has data rows 2-6 (row 1 = header). We are inserting a new
# blank row at row 3 (pushing the existing rows 3-6 down to rows 4-7), then
# filling in a missing "finish" time + total for the (now) row 7 entry, and
# appending a brand new row 8 for a new time-log entry.
#
# Work from the bottom up so we never overwrite source data before it is moved,
# copying formatting explicitly (via PasteSpecial of formats) so cells land on
# the same shared style indexes as their donor cells instead of minting new
# near-duplicate styles.

# --- old row 6 ("setting up github", 8/5/2018) becomes new row 7 ---
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A7").Value = $ws.Range("A6").Value()

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B7").Value = $ws.Range("B6").Value()

$ws.Range("D7").Value = $ws.Range("D6").Value()

# --- old row 5 (8/3/2018, locus questions) becomes new row 6 ---
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A6").Value = $ws.Range("A5").Value()

$ws.Range("B5").Copy() | Out-Null
$ws.Range("B6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B6").Value = $ws.Range("B5").Value()

$ws.Range("C5").Copy() | Out-Null
$ws.Range("C6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C6").Value = $ws.Range("C5").Value()

$ws.Range("D6").Value = $ws.Range("D5").Value()

$ws.Range("E5").Copy() | Out-Null
$ws.Range("E6").PasteSpecial($xlPasteFormats) | Out-Null

# --- old row 4 (no date, writing set 2 cont'd) becomes new row 5 ---
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B5").Value = $ws.Range("B4").Value()

$ws.Range("C4").Copy() | Out-Null
$ws.Range("C5").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C5").Value = $ws.Range("C4").Value()

$ws.Range("D5").Value = $ws.Range("D4").Value()

$ws.Range("E4").Copy() | Out-Null
$ws.Range("E5").PasteSpecial($xlPasteFormats) | Out-Null

# old row 5 had a date (A5) but new row 5 (old row 4) must not - clear it out
# now that its former content has already been relocated to A6 above.
$ws.Range("A5").Clear()

# --- old row 3 (7/30/2018, writing set 2) becomes new row 4 ---
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").Value = $ws.Range("A3").Value()

$ws.Range("B3").Copy() | Out-Null
$ws.Range("B4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B4").Value = $ws.Range("B3").Value()

$ws.Range("C3").Copy() | Out-Null
$ws.Range("C4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C4").Value = $ws.Range("C3").Value()

$ws.Range("D4").Value = $ws.Range("D3").Value()

$ws.Range("E3").Copy() | Out-Null
$ws.Range("E4").PasteSpecial($xlPasteFormats) | Out-Null

# Row 3 is now entirely stale (its content has been relocated to row 4) -
# wipe it completely so it no longer exists as a populated row.
$ws.Range("A3:G3").Clear()

# --- fill in the missing finish time for row 7 ("setting up github") ---
$ws.Range("B7").Copy() | Out-Null
$ws.Range("C7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C7").Value = 0.57291666666666663

# Re-establish the total-hours formula as a single multi-cell assignment so it
# is stored as one shared formula, matching the original file's E3:E5 shared
# group (now shifted down to E4:E6).
$ws.Range("E4:E6").Formula = "=C4-B4"

# Row 7's total is a brand new, individually-typed formula (not part of the
# shared-formula fill), same treatment as row 8 below.
$ws.Range("B7").Copy() | Out-Null
$ws.Range("E7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E7").Formula = "=C7-B7"

# --- brand new row 8: "worked solution; commute q." ---
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("B8").Value = 0.58333333333333337

$ws.Range("B7").Copy() | Out-Null
$ws.Range("C8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C8").Value = 0.71875

$ws.Range("D8").Value = "worked solution; commute q."

$ws.Range("B7").Copy() | Out-Null
$ws.Range("E8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E8").Formula = "=C8-B8"

# Widen column D to fit the new, longer description text.
# (The runtime's ColumnWidth -> stored OOXML width conversion adds ~5/6 of a
# character, so back that off here to land on a stored width of exactly 30.)
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668

# Match the author's final selection.
$ws.Range("C9").Select()
